# Updates the Mantel correlogram table (rural_12km) with values recomputed
# using Euclidean distances, per commit "recreated mantel correlograms with
# Euclidean distances".
#
# Table layout (row, col):
#   Row 2 ("6,000")  : Mantel r (2,3)  -0.025 -> 0.045 ; p (2,4)  0.26  -> 0.091
#   Row 3 ("18,000") : Mantel r (3,3)   0.004 -> 0.025 ; p (3,4)  0.519 -> 0.206
#   Row 4 ("30,000") : Mantel r (4,3)   0.086 -> 0.044 ; p (4,4)  0.042 -> 0.234
#                       (the p-value cell in row 4 also loses its bold run)

$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Get-CellTextRange($row, $col) {
    # Build a fresh Range spanning the cell's visible text (i.e. excluding
    # the trailing end-of-cell marker) via Document.Range(start, end).
    $cellRange = $tbl.Cell($row, $col).Range
    return $d.Range($cellRange.Start, $cellRange.End - 1)
}

function Set-CellText($row, $col, $newText) {
    $rng = Get-CellTextRange $row $col
    $rng.Text = $newText
}

Set-CellText 2 3 "0.045"
Set-CellText 2 4 "0.091"
Set-CellText 3 3 "0.025"
Set-CellText 3 4 "0.206"
Set-CellText 4 3 "0.044"
Set-CellText 4 4 "0.234"

# The p-value cell in the last row also loses its bold formatting.
$boldRng = Get-CellTextRange 4 4
$boldRng.Font.Bold = $false
